# Add "2022-Q3" data to the workbook:
#  1. Insert a brand-new worksheet named "2022-Q3" right after "总计" (i.e. before the current
#     "2022-Q2" sheet). All the other quarter tabs shift right by one position automatically,
#     carrying their existing data with them unchanged.
#  2. Populate the new "2022-Q3" sheet with the fund holdings for that quarter.
#  3. Update the "总计" (summary) sheet: insert a new row for "2022-Q3" at the top of the data
#     (row 2), pushing the other rows down, and recompute the sequential index column (A).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: insert the new "2022-Q3" worksheet before the existing second sheet ("2022-Q2").
# ---------------------------------------------------------------------------
$beforeSheet = $wb.Worksheets.Item(2)
$q3 = $wb.Worksheets.Add($beforeSheet)
$q3.Name = "2022-Q3"

# ---------------------------------------------------------------------------
# Step 2: fill in the "2022-Q3" fund table.
# ---------------------------------------------------------------------------

# Match the look of the other fund-holding sheets: copy the header/index-column
# formatting from the "2022-Q2" sheet (its header row + column A use the same style).
$styleSrc = $wb.Worksheets.Item("2022-Q2")
$styleSrc.Range("B1:H1").Copy($q3.Range("B1:H1"))
$styleSrc.Range("A2").Copy($q3.Range("A2:A9"))

$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

$q3Rows = @(
    @{A=0; B="166005"; C="中欧价值发现混合 -A";            D="26.62"; E="93.73"; F="4.86"; G="1.2937"; H=7},
    @{A=1; B="001810"; C="中欧潜力价值灵活配置混合A";        D="19.07"; E="93.66"; F="4.95"; G="0.9440"; H=5},
    @{A=2; B="004232"; C="中欧价值发现混合 -C";            D="8.18";  E="93.73"; F="4.86"; G="0.3975"; H=7},
    @{A=3; B="166024"; C="中欧恒利三年定期开放混合";         D="3.99";  E="98.45"; F="4.79"; G="0.1911"; H=4},
    @{A=4; B="166020"; C="中欧成长优选回报灵活配置混合A";     D="2.38";  E="93.70"; F="4.67"; G="0.1111"; H=3},
    @{A=5; B="005764"; C="中欧潜力价值灵活配置混合C";        D="2.01";  E="93.66"; F="4.95"; G="0.0995"; H=5},
    @{A=6; B="001891"; C="中欧成长优选回报灵活配置混合E";     D="0.74";  E="93.70"; F="4.67"; G="0.0346"; H=3},
    @{A=7; B="001882"; C="中欧价值发现混合 -E";            D="0.43";  E="93.73"; F="4.86"; G="0.0209"; H=7}
)

$r = 2
foreach ($row in $q3Rows) {
    $q3.Cells.Item($r, 1).Value = $row.A

    $q3.Cells.Item($r, 2).NumberFormat = "@"
    $q3.Cells.Item($r, 2).Value = $row.B

    $q3.Cells.Item($r, 3).Value = $row.C

    $q3.Cells.Item($r, 4).NumberFormat = "@"
    $q3.Cells.Item($r, 4).Value = $row.D

    $q3.Cells.Item($r, 5).NumberFormat = "@"
    $q3.Cells.Item($r, 5).Value = $row.E

    $q3.Cells.Item($r, 6).NumberFormat = "@"
    $q3.Cells.Item($r, 6).Value = $row.F

    $q3.Cells.Item($r, 7).NumberFormat = "@"
    $q3.Cells.Item($r, 7).Value = $row.G

    $q3.Cells.Item($r, 8).Value = $row.H

    $r = $r + 1
}

# ---------------------------------------------------------------------------
# Step 3: update the "总计" summary sheet with the new 2022-Q3 row, shifting the rest down.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

# Row 9 is brand new; give its index cell (column A) the same style as the existing
# index cells above it (A2:A8 all share one style).
$summary.Range("A8").Copy($summary.Range("A9"))

$summaryRows = @(
    @{A=0; B="2022-Q3"; C=8;  D=3.09},
    @{A=1; B="2022-Q2"; C=8;  D=3.6},
    @{A=2; B="2022-Q1"; C=10; D=5.69},
    @{A=3; B="2021-Q4"; C=8;  D=2.63},
    @{A=4; B="2021-Q3"; C=8;  D=2.44},
    @{A=5; B="2021-Q2"; C=3;  D=0.23},
    @{A=6; B="2021-Q1"; C=1;  D=0.73},
    @{A=7; B="2020-Q4"; C=5;  D=1.26}
)

$r = 2
foreach ($row in $summaryRows) {
    $summary.Cells.Item($r, 1).Value = $row.A
    $summary.Cells.Item($r, 2).Value = $row.B
    $summary.Cells.Item($r, 3).Value = $row.C
    $summary.Cells.Item($r, 4).Value = $row.D
    $r = $r + 1
}
